# Updates the cryptos list worksheet (row 2..51) with refreshed price /
# 1h-volume-change figures, and swaps the WrappedEther / ShibaInu rows
# (16 <-> 17) back the other way around, per the upstream GitHub Actions
# refresh commit "Updated cryptos list on Sat Jun 29 12:45:12 UTC 2024".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All of column D (Price) on the data rows is stored as literal text in the
# workbook (e.g. "61.074.78", "7.60") rather than as numbers, so values that
# look numeric must be forced to Text first -- otherwise COM's Range.Value
# setter re-interprets them (dropping the "573.46" -> 573.46 style, or worse,
# turning "1.00" into 1). Stamp the whole Price column as Text, write the new
# values, then restore the "Normal" style so the cells end up indistinguishable
# from the untouched ones (no stray number-format override left behind).
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$updates = @{
    "D2" = '61.077.03'
    "E2" = '  -0.92%  '
    "D3" = '3.396.37'
    "E3" = '  -1.79%  '
    "E4" = '  +0.01%  '
    "D5" = '573.46'
    "E5" = '  -0.81%  '
    "D6" = '142.51'
    "E6" = '  -2.16%  '
    "D7" = '3.397.00'
    "E7" = '  -1.80%  '
    "E8" = '  +0.02%  '
    "E9" = '  -0.83%  '
    "D10" = '7.62'
    "E10" = '  -0.19%  '
    "E11" = '  -2.89%  '
    "E12" = '  +1.17%  '
    "D13" = '3.978.36'
    "E13" = '  -1.80%  '
    "E14" = '  +2.04%  '
    "D15" = '27.98'
    "E15" = '  -2.95%  '
    "B16" = 'ShibaInu'
    "C16" = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
    "D16" = '0.0000171'
    "E16" = '  -2.20%  '
    "B17" = 'WrappedEther'
    "C17" = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
    "D17" = '3.395.41'
    "E17" = '  -2.27%  '
    "D18" = '61.106.83'
    "E18" = '  -1.01%  '
    "E19" = '  -4.00%  '
    "D20" = '13.81'
    "E20" = '  -4.01%  '
    "D21" = '8.95'
    "E21" = '  -5.27%  '
    "D22" = '382.99'
    "E22" = '  -4.89%  '
    "D23" = '0.557'
    "E23" = '  -1.81%  '
    "D24" = '74.43'
    "E24" = '  +0.49%  '
    "D25" = '1.00'
    "E25" = '  +0.45%  '
    "E26" = '  -5.82%  '
    "D27" = '3.536.29'
    "E27" = '  -1.65%  '
    "E28" = '  -0.50%  '
    "E29" = '  -0.16%  '
    "E30" = '  -3.74%  '
    "E31" = '  -3.30%  '
    "E32" = '  -1.72%  '
    "E33" = '  -6.15%  '
    "D35" = '23.45'
    "E35" = '  -2.49%  '
    "D36" = '6.99'
    "E36" = '  -1.68%  '
    "D37" = '167.23'
    "E37" = '  +0.02%  '
    "D38" = '3.429.61'
    "E38" = '  -1.60%  '
    "D39" = '4.99'
    "E39" = '  -4.06%  '
    "E40" = '  -5.24%  '
    "D41" = '0.0771'
    "E41" = '  -3.37%  '
    "D42" = '27.35'
    "E42" = '  -0.55%  '
    "D43" = '0.782'
    "E43" = '  -2.99%  '
    "E44" = '  +0.19%  '
    "D45" = '4.42'
    "E45" = '  -2.55%  '
    "D46" = '1.68'
    "E46" = '  -4.17%  '
    "D47" = '1.13'
    "E47" = '  -2.16%  '
    "D48" = '2.488.47'
    "E48" = '  -4.93%  '
    "D49" = '6.81'
    "E49" = '  -2.63%  '
    "D50" = '22.96'
    "E50" = '  -1.19%  '
    "D51" = '0.0265'
    "E51" = '  -0.06%  '
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

$priceRange.Style = "Normal"
